# fix de formato pptx
# - Tweaks the column widths of the "Unidad de medida" tables on both
#   slides (first column narrower, second column wider).
# - Nudges the "Ficha Técnica" title textbox on slide 2 a touch to the
#   left/down.

$p = $ppt.ActivePresentation

# --- Slide 1 : "Tabla 16" ---------------------------------------------
$s1 = $p.Slides.Item(1)
$tbl1 = $s1.Shapes.Item(2).Table
$tbl1.Columns.Item(1).Width = 111.7172440944882
$tbl1.Columns.Item(2).Width = 373.24826771653545

# --- Slide 2 : "Tabla 12" ---------------------------------------------
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(2).Table
$tbl2.Columns.Item(1).Width = 111.7172440944882
$tbl2.Columns.Item(2).Width = 373.24826771653545

# --- Slide 2 : "CuadroTexto 16" (Ficha Técnica title) -------------------
$title = $s2.Shapes.Item(7)
$title.Left = 817.5085826771654
$title.Top = 54.749212598425196
